# Apply the "final evaluation form" update to the Computer Graphics
# Homework 1 evaluation form workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluation form")

# --- H column remarks: append extra text to existing remarks ---

$ws.Range("H21").Value = "Circle geometry, box geometry, cylinder geometry, sphere geometry, torus geometry, custom angle geometry, shape geometry"

$ws.Range("H23").Value = "The tree has 3 different colors applied, lighter-brown for the tree trunk, darker brown for the tree ground and green for the tree's leafs (see: src/world/environment/tree.js). The floor and road make use of textures (also material properties) which comes down to a total of 5. The final material property is applied for the finish line,  which is a white color for the white blocks. Total material properties: 6."
$ws.Rows.Item(23).RowHeight = 165

$ws.Range("H24").Value = "We have a texture for the floor, and a texture for the road. (skybox kind of counts)"
$ws.Rows.Item(24).RowHeight = 45

$ws.Range("H26").Value = "The car is animated (movement, and wheels), and the lights (light goes from red to green)"
$ws.Rows.Item(26).RowHeight = 45

$ws.Range("H33").Value = "lots of trees :)"

$ws.Range("H34").Value = "We think that a racing track with a moving car is pretty cool. Look in the race car when it isn't moving to see even more coolness :). What is also cool: The light of the sun is aligned with the sun of the skybox. The car's driver (Eltjo) also respects the traffic laws(stops at red light, continues at green light)!"
$ws.Rows.Item(34).RowHeight = 120.75

# --- Row 6: Student 2 number gets filled in ---
$ws.Range("E6").Value = "s1127251"
$ws.Rows.Item(6).RowHeight = 30

# --- Update the view: scroll position and active selection ---
$ws.Application.Goto($ws.Range("A26"))
$ws.Range("H34").Select()

$wb.Save()
